# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> linked only from the Notes Master (name="Office Theme")
#   ppt/theme/theme2.xml -> linked from the Slide Master + the presentation itself (name="Integral")
#
# The target edit swaps the two themes' colour schemes (the "Integral"/Red
# Violet palette becomes the deck's main theme colours and the original
# "Office"/blue-orange palette is what had been the secondary theme) while
# every other part of each theme (font scheme, format scheme, relationship
# wiring) is left untouched. Apply it the way PowerPoint itself would, by
# pushing the new colour values onto the Slide Master's ThemeColorScheme
# (the deck's active design) one swatch at a time.

function Set-ThemeColorRGB {
    param(
        $ColorScheme,
        [int]$Index,
        [string]$Hex
    )
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    $ColorScheme.Item($Index).RGB = $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
Set-ThemeColorRGB $colors 1  "000000"
Set-ThemeColorRGB $colors 2  "FFFFFF"
Set-ThemeColorRGB $colors 3  "44546A"
Set-ThemeColorRGB $colors 4  "E7E6E6"
Set-ThemeColorRGB $colors 5  "5B9BD5"
Set-ThemeColorRGB $colors 6  "ED7D31"
Set-ThemeColorRGB $colors 7  "A5A5A5"
Set-ThemeColorRGB $colors 8  "FFC000"
Set-ThemeColorRGB $colors 9  "4472C4"
Set-ThemeColorRGB $colors 10 "70AD47"
Set-ThemeColorRGB $colors 11 "0563C1"
Set-ThemeColorRGB $colors 12 "954F72"
